# Non-Oncology import tool scenarios
# Replace the old generic "Report-" filename-prefix values with fully
# qualified "StandardExcelReport-...-2023_" values on each scenario sheet,
# and move the active selection/tab to reflect where the author left off.

$wb = $excel.ActiveWorkbook

# --- NewImportLogic sheet (QOL_and_ECON - UtilityOutcome scenario) ---
$ws1 = $wb.Worksheets.Item("NewImportLogic")
$ws1.Range("H2").Value = "StandardExcelReport-QOL_and_ECON - UtilityOutcome-Quality of Life-2023_"

# --- OldImportLogic sheet (Alkermes - Melanoma scenario) ---
$ws2 = $wb.Worksheets.Item("OldImportLogic")
$ws2.Range("H2").Value = "StandardExcelReport-Alkermes - Melanoma-Quality of Life-2023_"

# --- prodfix sheet (PRODFix_QOL_ECON - UtilityOutcome scenario) ---
$ws3 = $wb.Worksheets.Item("prodfix")
$ws3.Range("L4").Value = "StandardExcelReport-PRODFix_QOL_ECON - UtilityOutcome-Quality of Life-2023_"

# --- Update each sheet's selection, then finish with prodfix active/selected ---
$ws1.Activate()
$ws1.Range("H2").Select()

$ws2.Activate()
$ws2.Range("H2").Select()

$ws3.Activate()
$ws3.Range("L6").Select()
